$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.996049335404507
$ws.Range("B3").Value = 0.9012017778810559
$ws.Range("B4").Value = 0.8432139183692016
$ws.Range("B5").Value = 0.8196463723957379
$ws.Range("B6").Value = 0.8157368148755779
$ws.Range("B7").Value = 0.8428958226820669
$ws.Range("B8").Value = 0.9632944265005676
$ws.Range("B9").Value = 1.201367543761194
$ws.Range("B10").Value = 1.377497620503902
$ws.Range("B11").Value = 1.457893145323396
$ws.Range("B12").Value = 1.488376124515923
$ws.Range("B13").Value = 1.481809342697886
$ws.Range("B14").Value = 1.460400219771884
$ws.Range("B15").Value = 1.447291575443217
$ws.Range("B16").Value = 1.372249050133348
$ws.Range("B17").Value = 1.326282590008532
$ws.Range("B18").Value = 1.299869618627326
$ws.Range("B19").Value = 1.29093106830868
$ws.Range("B20").Value = 1.331173138610666
$ws.Range("B21").Value = 1.466687543852004
$ws.Range("B22").Value = 1.555480870918473
$ws.Range("B23").Value = 1.508069548555511
$ws.Range("B24").Value = 1.328962078709026
$ws.Range("B25").Value = 1.136749828668144
$ws.Range("C2").Value = 0.1829687567014275
$ws.Range("C3").Value = 0.1600504625861845
$ws.Range("C4").Value = 0.1459627490864364
$ws.Range("C5").Value = 0.1402180414634984
$ws.Range("C6").Value = 0.1392639082839366
$ws.Range("C7").Value = 0.1458852893734388
$ws.Range("C8").Value = 0.1750698690751449
$ws.Range("C9").Value = 0.2321733053422008
$ws.Range("C10").Value = 0.2740512108289295
$ws.Range("C11").Value = 0.2930868018417243
$ws.Range("C12").Value = 0.3002929185513779
$ws.Range("C13").Value = 0.298741054053977
$ws.Range("C14").Value = 0.2936796990651658
$ws.Range("C15").Value = 0.2905791770464532
$ws.Range("C16").Value = 0.2728068813363507
$ws.Range("C17").Value = 0.2619002633106788
$ws.Range("C18").Value = 0.2556256707910336
$ws.Range("C19").Value = 0.2535009655327372
$ws.Range("C20").Value = 0.2630614365266126
$ws.Range("C21").Value = 0.2951664035740578
$ws.Range("C22").Value = 0.3161357035396577
$ws.Range("C23").Value = 0.3049452312641279
$ws.Range("C24").Value = 0.2625364832940136
$ws.Range("C25").Value = 0.2167388039821105
$ws.Range("D2").Value = 0.1878679431105965
$ws.Range("D3").Value = 0.1866206637211718
$ws.Range("D4").Value = 0.1859009606257871
$ws.Range("D5").Value = 0.1856193327121929
$ws.Range("D6").Value = 0.1855732743477461
$ws.Range("D7").Value = 0.1858971152181255
$ws.Range("D8").Value = 0.1874283397786911
$ws.Range("D9").Value = 0.1907948997231728
$ws.Range("D10").Value = 0.1934875938347034
$ws.Range("D11").Value = 0.1947596920133208
$ws.Range("D12").Value = 0.195248138777373
$ws.Range("D13").Value = 0.1951426445137372
$ws.Range("D14").Value = 0.1947997422303018
$ws.Range("D15").Value = 0.1945905796591774
$ws.Range("D16").Value = 0.1934054039881659
$ws.Range("D17").Value = 0.192690382297755
$ws.Range("D18").Value = 0.1922835644952201
$ws.Range("D19").Value = 0.1921465878391828
$ws.Range("D20").Value = 0.1927660380530796
$ws.Range("D21").Value = 0.1949002786924297
$ws.Range("D22").Value = 0.1963343363985928
$ws.Range("D23").Value = 0.1955653826702815
$ws.Range("D24").Value = 0.1927318208211517
$ws.Range("D25").Value = 0.1898454690732407
$ws.Range("F2").Value = 1.765463493389589
$ws.Range("F3").Value = 1.757592780051567
$ws.Range("F4").Value = 1.753975854694346
$ws.Range("F5").Value = 1.752806690648669
$ws.Range("F6").Value = 1.7526309312602
$ws.Range("F7").Value = 1.753958854224948
$ws.Range("F8").Value = 1.762496773407918
$ws.Range("F9").Value = 1.788931703733311
$ws.Range("F10").Value = 1.814330436216977
$ws.Range("F11").Value = 1.827197921082387
$ws.Range("F12").Value = 1.832260514456195
$ws.Range("F13").Value = 1.831161728786725
$ws.Range("F14").Value = 1.827610610478573
$ws.Range("F15").Value = 1.82546021801204
$ws.Range("F16").Value = 1.813516038497994
$ws.Range("F17").Value = 1.806525850451777
$ws.Range("F18").Value = 1.802628844740028
$ws.Range("F19").Value = 1.801330575224455
$ws.Range("F20").Value = 1.807257171558291
$ws.Range("F21").Value = 1.828648496018076
$ws.Range("F22").Value = 1.843736677643037
$ws.Range("F23").Value = 1.835582117017239
$ws.Range("F24").Value = 1.806926162138936
$ws.Range("F25").Value = 1.780735691266244
$ws.Range("G2").Value = 1.093982151336363
$ws.Range("G3").Value = 1.083478378094483
$ws.Range("G4").Value = 1.077945592247332
$ws.Range("G5").Value = 1.075920273241223
$ws.Range("G6").Value = 1.075597788454772
$ws.Range("G7").Value = 1.077917351034927
$ws.Range("G8").Value = 1.090169520924206
$ws.Range("G9").Value = 1.12152293348305
$ws.Range("G10").Value = 1.149104206321567
$ws.Range("G11").Value = 1.162656067776538
$ws.Range("G12").Value = 1.167933649102338
$ws.Range("G13").Value = 1.166790523464812
$ws.Range("G14").Value = 1.163087328537841
$ws.Range("G15").Value = 1.160838037968375
$ws.Range("G16").Value = 1.148238891781233
$ws.Range("G17").Value = 1.14076802954591
$ws.Range("G18").Value = 1.136565506798746
$ws.Range("G19").Value = 1.135158800067671
$ws.Range("G20").Value = 1.141553524014
$ws.Range("G21").Value = 1.164171079959374
$ws.Range("G22").Value = 1.17980330051634
$ws.Range("G23").Value = 1.171381864318846
$ws.Range("G24").Value = 1.141198113699204
$ws.Range("G25").Value = 1.112248426180912
$ws.Range("H2").Value = 1.054649141061986
$ws.Range("H3").Value = 1.056222366358327
$ws.Range("H4").Value = 1.05784645235471
$ws.Range("H5").Value = 1.058673376029972
$ws.Range("H6").Value = 1.058820645930936
$ws.Range("H7").Value = 1.057856936602221
$ws.Range("H8").Value = 1.055054756652325
$ws.Range("H9").Value = 1.054800549551487
$ws.Range("H10").Value = 1.057836773893229
$ws.Range("H11").Value = 1.059924053698467
$ws.Range("H12").Value = 1.060816462066015
$ws.Range("H13").Value = 1.060619721985603
$ws.Range("H14").Value = 1.059995425743438
$ws.Range("H15").Value = 1.059626323671722
$ws.Range("H16").Value = 1.057714613325118
$ws.Range("H17").Value = 1.056723017644117
$ws.Range("H18").Value = 1.05621911018531
$ws.Range("H19").Value = 1.056059891381494
$ws.Range("H20").Value = 1.056821695874191
$ws.Range("H21").Value = 1.060176024679322
$ws.Range("H22").Value = 1.062963029560564
$ws.Range("H23").Value = 1.061420977500205
$ws.Range("H24").Value = 1.056776877352746
$ws.Range("H25").Value = 1.054305360844637
$ws.Range("J2").Value = 0.2575401434904663
$ws.Range("J3").Value = 0.2577735803831658
$ws.Range("J4").Value = 0.2580473444133418
$ws.Range("J5").Value = 0.2581916641827746
$ws.Range("J6").Value = 0.2582176057653456
$ws.Range("J7").Value = 0.2580491581737192
$ws.Range("J8").Value = 0.2575935366278728
$ws.Range("J9").Value = 0.2577370685218838
$ws.Range("J10").Value = 0.2584780024960978
$ws.Range("J11").Value = 0.2589538081555887
$ws.Range("J12").Value = 0.2591539902977402
$ws.Range("J13").Value = 0.2591099869248197
$ws.Range("J14").Value = 0.2589698760299441
$ws.Range("J15").Value = 0.2588866608578471
$ws.Range("J16").Value = 0.258449703987587
$ws.Range("J17").Value = 0.2582172198235853
$ws.Range("J18").Value = 0.2580965578075691
$ws.Range("J19").Value = 0.2580579446331654
$ws.Range("J20").Value = 0.2582406164603128
$ws.Range("J21").Value = 0.2590104866471421
$ws.Range("J22").Value = 0.2596302698741724
$ws.Range("J23").Value = 0.2592887908213442
$ws.Range("J24").Value = 0.2582299983639373
$ws.Range("J25").Value = 0.257586882418984
$ws.Range("M2").Value = 0.4394676603794849
$ws.Range("M3").Value = 0.4143720984839732
$ws.Range("M4").Value = 0.3991453431329219
$ws.Range("M5").Value = 0.3929861922182809
$ws.Range("M6").Value = 0.39196624391478
$ws.Range("M7").Value = 0.3990620927772
$ws.Range("M8").Value = 0.4307769845183032
$ws.Range("M9").Value = 0.4944133442555199
$ws.Range("M10").Value = 0.5420518457471104
$ws.Range("M11").Value = 0.5639174779473706
$ws.Range("M12").Value = 0.5722254263587132
$ws.Range("M13").Value = 0.5704349210151634
$ws.Range("M14").Value = 0.5646004189825504
$ws.Range("M15").Value = 0.5610302501912372
$ws.Range("M16").Value = 0.5406267848281061
$ws.Range("M17").Value = 0.5281597117199013
$ws.Range("M18").Value = 0.5210073057445541
$ws.Range("M19").Value = 0.5185887738424313
$ws.Range("M20").Value = 0.5294849570006832
$ws.Range("M21").Value = 0.56631339685336
$ws.Range("M22").Value = 0.590545548138607
$ws.Range("M23").Value = 0.577597536560404
$ws.Range("M24").Value = 0.5288857666314257
$ws.Range("M25").Value = 0.4770430060050401
